$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

# Row 2
Set-TextCell "D2" "43.202.59"
Set-TextCell "E2" "  +2.42%  "

# Row 3
Set-TextCell "D3" "2.307.10"
Set-TextCell "E3" "  +1.59%  "

# Row 4
Set-TextCell "E4" "  +0.07%  "

# Row 5
Set-TextCell "D5" "302.07"
Set-TextCell "E5" "  +1.32%  "

# Row 6
Set-TextCell "E6" "  +5.77%  "

# Row 7
Set-TextCell "E7" "  +2.51%  "

# Row 8
Set-TextCell "E8" "  -0.02%  "

# Row 9
Set-TextCell "E9" "  +3.37%  "

# Row 10
Set-TextCell "D10" "34.43"
Set-TextCell "E10" "  +4.07%  "

# Row 11
Set-TextCell "E11" "  +1.19%  "

# Row 12
Set-TextCell "E12" "  +4.47%  "

# Row 13
Set-TextCell "D13" "17.99"
Set-TextCell "E13" "  +14.66%  "

# Row 14
Set-TextCell "E14" "  +2.72%  "

# Row 15
Set-TextCell "D15" "2.666.17"

# Row 16
Set-TextCell "D16" "2.335.35"
Set-TextCell "E16" "  +3.79%  "

# Row 17
Set-TextCell "E17" "  +5.57%  "

# Row 18
Set-TextCell "D18" "43.117.47"
Set-TextCell "E18" "  +2.25%  "

# Row 19
Set-TextCell "D19" "12.56"
Set-TextCell "E19" "  +10.76%  "

# Row 20
Set-TextCell "D20" "0.0₃0907"
Set-TextCell "E20" "  +2.29%  "

# Row 21
Set-TextCell "E21" "  +2.49%  "

# Row 22
Set-TextCell "D22" "67.80"
Set-TextCell "E22" "  +1.84%  "

# Row 23
Set-TextCell "D23" "237.18"
Set-TextCell "E23" "  +1.91%  "

# Row 24
Set-TextCell "E24" "  +14.77%  "

# Row 25
Set-TextCell "E25" "  +0.34%  "

# Row 26
Set-TextCell "E26" "  +0.80%  "

# Row 27
Set-TextCell "D27" "24.84"
Set-TextCell "E27" "  +4.46%  "

# Row 28
Set-TextCell "D28" "168.76"
Set-TextCell "E28" "  +1.50%  "

# Row 29
Set-TextCell "E29" "  -9.35%  "

# Row 30
Set-TextCell "D30" "34.26"
Set-TextCell "E30" "  +1.64%  "

# Row 31
Set-TextCell "E31" "  +1.60%  "

# Row 32
Set-TextCell "E32" "  +0.12%  "

# Row 33
Set-TextCell "D33" "5.04"
Set-TextCell "E33" "  +2.82%  "

# Row 34
Set-TextCell "E34" "  +4.42%  "

# Row 35
Set-TextCell "D35" "4.54"
Set-TextCell "E35" "  +1.41%  "

# Row 36
Set-TextCell "D36" "17.08"
Set-TextCell "E36" "  +6.57%  "

# Row 37
Set-TextCell "D37" "0.0691"
Set-TextCell "E37" "  +0.52%  "

# Row 38
Set-TextCell "E38" "  +3.93%  "

# Row 39
Set-TextCell "E39" "  +5.47%  "

# Row 40
Set-TextCell "D40" "2.82"
Set-TextCell "E40" "  +1.75%  "

# Row 41
Set-TextCell "E41" "  +0.79%  "

# Row 42
Set-TextCell "E42" "  -3.49%  "

# Row 43
Set-TextCell "D43" "2.001.77"
Set-TextCell "E43" "  +2.23%  "

# Row 45
Set-TextCell "D45" "10.15"
Set-TextCell "E45" "  +6.40%  "

# Row 46
Set-TextCell "D46" "17.63"
Set-TextCell "E46" "  +1.56%  "

# Row 47
Set-TextCell "E47" "  +2.81%  "

# Row 48
Set-TextCell "D48" "56.11"
Set-TextCell "E48" "  +8.29%  "

# Row 49
Set-TextCell "D49" "2.534.21"
Set-TextCell "E49" "  +1.55%  "

# Row 50
Set-TextCell "E50" "  +5.57%  "

# Row 51
Set-TextCell "E51" "  +1.54%  "
